# POZNAMKY/55. CYKLUS-FOR-EACH.docx
# -----------------------------------------------------------------------
# The document's title heading currently reads just "CYKLUS forEach".
# This note is lesson 55, so prefix the heading with "55. " using the
# same run formatting (Consolas, bold, 28pt/14pt, single underline,
# Arial for complex scripts) as the rest of the heading run.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# Locate the heading text and collapse the found range to its start,
# right before the "C" of "CYKLUS forEach".
$heading = $d.Content
$heading.Find.Execute("CYKLUS forEach", $false, $false, $false, $false,
                       $false, $true, 1, $false, "", 0)
$heading.Collapse(1)

# Insert the new lesson-number text in front of the heading.
$heading.InsertBefore("55. ")

# Re-grab exactly the text we just inserted ("55. ") and nudge a
# formatting property off/on (to the same value it already has by
# inheritance) so the engine keeps it as its own run instead of
# silently merging it back into the following "CYKLUS forEach" run.
$start = $heading.Start
$newRun = $d.Range($start, $start + 4)
$newRun.Font.Bold = $false
$newRun.Font.Bold = $true
